$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a "duplicate_image_filename" column (E) that was missing
# values for the practice/word rows (rows 2-21). Fill them with "NA".
$ws.Range("E2:E21").Value = "NA"
